# إضافة حدث جديد في Card7 by admin at 2026-01-20 16:05:08
#
# Net content change:
#  - Card7: a new service-visit event row (row 14) is appended:
#        Date=20/1/2026, Event="زياره توكيل",
#        Correction="تم تغير سوفت كرد لbc", Serviced by="م. احمد علي توكيل"
#    Previously-blank detail cells in rows 2-13 are re-stamped with the
#    literal text "nan" (artifact of the sync tool that also wrote this row).
#  - Card6: the same event row had been mistakenly duplicated there
#    (present on both its row 13 and row 14). The duplicate row 14 is
#    removed and the stray "nan" filler text that came with it is cleared
#    back to blank, restoring Card6 to its normal 13-row shape.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Card7: stamp previously-empty cells with "nan" and append the new row
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Card7")

$row2Cols  = @("D","E","F","G","H","I","J","K","N")
$row3Cols  = @("G","H","I","J","K")
$fullCols  = @("D","E","F","G","H","I","J","K","L","M","N","O")
$row13Cols = @("B","C","D","E","F","G","H","I","J","K")

foreach ($col in $row2Cols) {
    $ws7.Range($col + "2").Value = "nan"
}

foreach ($col in $row3Cols) {
    $ws7.Range($col + "3").Value = "nan"
}

for ($r = 4; $r -le 12; $r++) {
    foreach ($col in $fullCols) {
        $ws7.Range($col + $r).Value = "nan"
    }
}

foreach ($col in $row13Cols) {
    $ws7.Range($col + "13").Value = "nan"
}

# New row 14 — the newly-logged event
$ws7.Range("A14").Value = "'7"
$ws7.Range("B14").Value = ""
$ws7.Range("C14").Value = ""
$ws7.Range("D14").Value = ""
$ws7.Range("E14").Value = ""
$ws7.Range("F14").Value = ""
$ws7.Range("G14").Value = ""
$ws7.Range("H14").Value = ""
$ws7.Range("I14").Value = ""
$ws7.Range("J14").Value = ""
$ws7.Range("K14").Value = ""
$ws7.Range("L14").Value = "20/1/2026"
$ws7.Range("M14").Value = "زياره توكيل"
$ws7.Range("N14").Value = "تم تغير سوفت كرد لbc"
$ws7.Range("O14").Value = "م. احمد علي توكيل"

# ---------------------------------------------------------------------
# Card6: remove the duplicate row and the stray "nan" filler text
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Card6")

foreach ($col in $fullCols) {
    $ws6.Range($col + "2").Value = ""
}

foreach ($col in $row3Cols) {
    $ws6.Range($col + "3").Value = ""
}

for ($r = 4; $r -le 12; $r++) {
    foreach ($col in $fullCols) {
        $ws6.Range($col + $r).Value = ""
    }
}

foreach ($col in $row13Cols) {
    $ws6.Range($col + "13").Value = ""
}

# Drop the duplicated event row entirely
$ws6.Rows.Item(14).Delete()
